$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.420.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.98%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.566.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'208.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.502"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'22.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.88%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.68%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0592"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.789.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.573.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.92%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.74%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'63.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.48%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'27.425.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.93%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'213.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.97%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0689"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.77%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.60%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +2.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.55%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'14.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.99%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.51%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.82%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D33").Value = "'1.376.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.66%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.28%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +1.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.61%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.532"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.21%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.21%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'63.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.39%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.701.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'85.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0₇0989"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0959"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.90%  "
$ws.Range("E51").Style = "Normal"
